# Updated remaining queries for C3DC
# Replaces the stale "id"-based JOIN conditions in every embedded SQL query
# with the corrected "study_id" / "participant_id" based JOIN conditions,
# resizes column C, and updates the sheet's active selection/scroll state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old JOIN block shared by every query in the sheet (Excel/COM always
# normalizes embedded line breaks to `n when read through Value2/Formula).
$oldJoin = "LEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"" + `
    "`nLEFT JOIN `n    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"" + `
    "`nLEFT JOIN `n    df_treatments trt ON prt.id = trt.`"participant.id`"" + `
    "`nLEFT JOIN `n    df_treatment_resp trr ON prt.id = trr.`"participant.id`"" + `
    "`nLEFT JOIN `n    df_survival srv ON prt.id = srv.`"participant.id`"" + `
    "`nLEFT JOIN `n    df_reference_files rfs ON std.id = rfs.`"study.id`""

# The corrected JOIN block.
$newJoin = "LEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"" + `
    "`nLEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"" + `
    "`nLEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"" + `
    "`nLEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"" + `
    "`nLEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"" + `
    "`nLEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`""

# Every cell that holds one of the seven SQL queries referencing the old JOIN.
$queryCells = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $queryCells) {
    $rng = $ws.Range($addr)
    $current = $rng.Value2
    if ($current -ne $null -and $current.Contains($oldJoin)) {
        $rng.Value2 = $current.Replace($oldJoin, $newJoin)
    }
}

# Column C no longer needs to "best fit" - it now has an explicit width.
$ws.Columns.Item(3).ColumnWidth = 68.16666666666667

# Move the active selection / scroll position down to row 6-7, matching the
# updated sheet view saved with the workbook.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("C7").Select()
